$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.551.68"
$ws.Range("E2").Value = "  -2.46%  "

$ws.Range("D3").Value = "'1.813.72"
$ws.Range("E3").Value = "  -2.20%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.77%  "

$ws.Range("D5").Value = "'1.008"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'308.82"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D7").Value = "'0.4573"
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("D8").Value = "'0.3666"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").Value = "'0.07138"
$ws.Range("E9").Value = "  -2.18%  "

$ws.Range("D10").Value = "'0.8803"
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").Value = "'0.07757"
$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("D12").Value = "'19.40"
$ws.Range("E12").Value = "  -3.32%  "

$ws.Range("D13").Value = "'1.824.51"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "'5.295"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "'6.378"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").Value = "'86.63"
$ws.Range("E16").Value = "  -5.13%  "

$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "'0.000008598"
$ws.Range("E18").Value = "  -3.56%  "

$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").Value = "'26.609.87"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").Value = "'5.012"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("D23").Value = "'10.48"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").Value = "'1.988"
$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("D25").Value = "'151.53"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'17.94"
$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("D27").Value = "'2.063"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("D28").Value = "'113.03"
$ws.Range("E28").Value = "  -2.39%  "

$ws.Range("D29").Value = "'4.848"
$ws.Range("E29").Value = "  -3.69%  "

$ws.Range("D30").Value = "'0.08695"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").Value = "'3.037"
$ws.Range("E31").Value = "  -3.30%  "

$ws.Range("D32").Value = "'4.521"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").Value = "'0.7341"
$ws.Range("E33").Value = "  -4.51%  "

$ws.Range("D34").Value = "'2.699"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "  -3.99%  "

$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -1.97%  "

$ws.Range("D38").Value = "'0.01955"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("D39").Value = "'0.05113"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").Value = "'6.984"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").Value = "'0.4992"
$ws.Range("E42").Value = "  -2.26%  "

$ws.Range("D43").Value = "'0.1558"
$ws.Range("E43").Value = "  -4.16%  "

$ws.Range("D44").Value = "'8.167"
$ws.Range("E44").Value = "  -3.90%  "

$ws.Range("D45").Value = "'1.009"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "'0.4603"
$ws.Range("E46").Value = "  -4.06%  "

$ws.Range("E47").Value = "  -3.52%  "

$ws.Range("D48").Value = "'101.11"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("D49").Value = "'1.590"
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("D50").Value = "'0.06000"
$ws.Range("E50").Value = "  -3.13%  "

$ws.Range("D51").Value = "'64.43"
$ws.Range("E51").Value = "  -1.30%  "
